# Applies the diff: adds a new "Sheet1" worksheet between "ReviewPaper" and
# "Colleges" containing a categorized index of links (General/ARIMA/SVM/ANN/LSTM),
# and appends a handful of matching "ANN" rows onto the "ReviewPaper" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add rows 19, 21, 22, 23 to the existing "ReviewPaper" sheet (ANN entry)
# ---------------------------------------------------------------------------
$review = $wb.Worksheets.Item("ReviewPaper")

$review.Range("A19").Value = "ANN"
$review.Rows.Item(19).RowHeight = 15.75

$review.Range("A21").Value = "https://www.sciencedirect.com/science/article/pii/S2314728817300715"
$review.Rows.Item(21).RowHeight = 15.75

$review.Range("A22").Value = "https://www.sciencedirect.com/science/article/pii/S1877050915006766"
$review.Rows.Item(22).RowHeight = 15.75

$review.Range("A23").Value = "https://www.researchgate.net/publication/318127667_Forecasting_of_nonlinear_time_series_using_ANN"
$review.Rows.Item(23).RowHeight = 15.75

# ---------------------------------------------------------------------------
# 2. Insert a new worksheet named "Sheet1" right after "ReviewPaper" (so that
#    the sheet order becomes ReviewPaper, Sheet1, Colleges)
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $review)
$newSheet.Name = "Sheet1"

# Category headers are typed first (matches original authoring order so that
# shared-string indices line up the same way the source workbook built them)
$newSheet.Range("A1").Value = "General"
$newSheet.Range("A6").Value = "ARIMA"
$newSheet.Range("A11").Value = "SVM"

# --- General ---------------------------------------------------------------
$newSheet.Range("A2").Value = "https://scholar.google.com/citations?user=vb9EOUMAAAAJ&hl=it"
$newSheet.Hyperlinks.Add($newSheet.Range("A2"), "https://scholar.google.com/citations?user=vb9EOUMAAAAJ&hl=it") | Out-Null
$newSheet.Range("A2").Style = "Hyperlink"

$newSheet.Range("A3").Value = "https://www.tableau.com/learn/articles/time-series-forecasting"
$newSheet.Hyperlinks.Add($newSheet.Range("A3"), "https://www.tableau.com/learn/articles/time-series-forecasting") | Out-Null
$newSheet.Range("A3").Style = "Hyperlink"

$newSheet.Range("A4").Value = "https://scholar.google.com/citations?view_op=view_citation&hl=it&user=vb9EOUMAAAAJ&citation_for_view=vb9EOUMAAAAJ:HeT0ZceujKMC"
$newSheet.Hyperlinks.Add($newSheet.Range("A4"), "https://scholar.google.com/citations?view_op=view_citation&hl=it&user=vb9EOUMAAAAJ&citation_for_view=vb9EOUMAAAAJ:HeT0ZceujKMC") | Out-Null
$newSheet.Range("A4").Style = "Hyperlink"

# --- ARIMA -------------------------------------------------------------
$newSheet.Range("A7").Value = "https://towardsdatascience.com/time-series-forecasting-predicting-stock-prices-using-an-arima-model-2e3b3080bd70"
$newSheet.Hyperlinks.Add($newSheet.Range("A7"), "https://towardsdatascience.com/time-series-forecasting-predicting-stock-prices-using-an-arima-model-2e3b3080bd70") | Out-Null
$newSheet.Range("A7").Style = "Hyperlink"

$newSheet.Range("A8").Value = "https://www.ijcsmc.com/docs/papers/August2016/V5I8201626.pdf"
$newSheet.Hyperlinks.Add($newSheet.Range("A8"), "https://www.ijcsmc.com/docs/papers/August2016/V5I8201626.pdf") | Out-Null
$newSheet.Range("A8").Style = "Hyperlink"

# --- SVM ---------------------------------------------------------------
$newSheet.Range("A12").Value = "https://towardsdatascience.com/lstm-time-series-forecasting-predicting-stock-prices-using-an-lstm-model-6223e9644a2f"
$newSheet.Hyperlinks.Add($newSheet.Range("A12"), "https://towardsdatascience.com/lstm-time-series-forecasting-predicting-stock-prices-using-an-lstm-model-6223e9644a2f") | Out-Null
$newSheet.Range("A12").Style = "Hyperlink"

$newSheet.Range("A13").Value = "https://www.researchgate.net/publication/222661136_Financial_time_series_forecasting_using_support_vector_machines"
$newSheet.Hyperlinks.Add($newSheet.Range("A13"), "https://www.researchgate.net/publication/222661136_Financial_time_series_forecasting_using_support_vector_machines") | Out-Null
$newSheet.Range("A13").Style = "Hyperlink"

$newSheet.Range("A14").Value = "https://sci-hub.mksa.top/10.1016/s0925-2312(03)00372-2"
$newSheet.Hyperlinks.Add($newSheet.Range("A14"), "https://sci-hub.mksa.top/10.1016/s0925-2312(03)00372-2") | Out-Null
$newSheet.Range("A14").Style = "Hyperlink"

$newSheet.Range("A15").Value = "https://www.cs.princeton.edu/sites/default/files/uploads/saahil_madge.pdf"
$newSheet.Hyperlinks.Add($newSheet.Range("A15"), "https://www.cs.princeton.edu/sites/default/files/uploads/saahil_madge.pdf") | Out-Null
$newSheet.Range("A15").Style = "Hyperlink"

$newSheet.Range("A16").Value = "https://www.analyticssteps.com/blogs/how-does-support-vector-machine-algorithm-works-machine-learning"
$newSheet.Hyperlinks.Add($newSheet.Range("A16"), "https://www.analyticssteps.com/blogs/how-does-support-vector-machine-algorithm-works-machine-learning") | Out-Null
$newSheet.Range("A16").Style = "Hyperlink"

$newSheet.Hyperlinks.Add($newSheet.Range("A17"), "https://jakevdp.github.io/PythonDataScienceHandbook/05.07-support-vector-machines.html", [System.Type]::Missing, [System.Type]::Missing, "https://jakevdp.github.io/PythonDataScienceHandbook/05.07-support-vector-machines.html") | Out-Null
$newSheet.Range("A17").Value = "In-Depth: Support Vector Machines | Python Data Science Handbook (jakevdp.github.io)"
$newSheet.Range("A17").Style = "Hyperlink"

# --- ANN -----------------------------------------------------------------
$newSheet.Range("A19").Value = "ANN"

$newSheet.Range("A20").Value = "https://jfin-swufe.springeropen.com/articles/10.1186/s40854-019-0131-7"
$newSheet.Hyperlinks.Add($newSheet.Range("A20"), "https://jfin-swufe.springeropen.com/articles/10.1186/s40854-019-0131-7") | Out-Null
$newSheet.Range("A20").Style = "Hyperlink"

$newSheet.Range("A21").Value = "https://www.sciencedirect.com/science/article/pii/S2314728817300715"
$newSheet.Hyperlinks.Add($newSheet.Range("A21"), "https://www.sciencedirect.com/science/article/pii/S2314728817300715") | Out-Null
$newSheet.Range("A21").Style = "Hyperlink"

$newSheet.Range("A22").Value = "https://www.sciencedirect.com/science/article/pii/S1877050915006766"
$newSheet.Hyperlinks.Add($newSheet.Range("A22"), "https://www.sciencedirect.com/science/article/pii/S1877050915006766") | Out-Null
$newSheet.Range("A22").Style = "Hyperlink"

$newSheet.Range("A23").Value = "https://www.researchgate.net/publication/318127667_Forecasting_of_nonlinear_time_series_using_ANN"
$newSheet.Hyperlinks.Add($newSheet.Range("A23"), "https://www.researchgate.net/publication/318127667_Forecasting_of_nonlinear_time_series_using_ANN") | Out-Null
$newSheet.Range("A23").Style = "Hyperlink"

# --- LSTM ------------------------------------------------------------------
$newSheet.Range("A25").Value = "LSTM"

$newSheet.Range("A26").Value = "https://www.researchgate.net/publication/348390803_Stock_Price_Prediction_Using_LSTM"
$newSheet.Hyperlinks.Add($newSheet.Range("A26"), "https://www.researchgate.net/publication/348390803_Stock_Price_Prediction_Using_LSTM") | Out-Null
$newSheet.Range("A26").Style = "Hyperlink"

$newSheet.Range("A27").Value = "https://towardsdatascience.com/lstm-time-series-forecasting-predicting-stock-prices-using-an-lstm-model-6223e9644a2f"
$newSheet.Hyperlinks.Add($newSheet.Range("A27"), "https://towardsdatascience.com/lstm-time-series-forecasting-predicting-stock-prices-using-an-lstm-model-6223e9644a2f") | Out-Null
$newSheet.Range("A27").Style = "Hyperlink"

# Row heights to match the default 12.45 look used throughout the new sheet
$newSheet.Rows.Item("1:27").RowHeight = 12.45

# View state: scrolled down a bit, focused on A10, selection on A29
$newSheet.Activate()
Write-Host "Done applying edits"
